$d = $word.ActiveDocument
$bm = $d.Bookmarks("_GoBack")
$bmRange = $bm.Range
$bmRange.InsertBefore(" Tentatvie gps + boussole sur epsom juste compass aui fonctionne")

# Re-split "30/05/2016" from what follows (restore original boundary)
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("30/05/2016", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Found1:" $found1 "Start/End:" $rng1.Start $rng1.End
$rng1.Bold = $true
$rng1.Bold = $false

# Split new text from what precedes it
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(" Tentatvie gps + boussole sur epsom juste compass aui fonctionne", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Found2:" $found2 "Start/End:" $rng2.Start $rng2.End
$rng2.Bold = $true
$rng2.Bold = $false
